$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title near the top of the document.
$d.Paragraphs(2).Range.Delete()

# 2. Insert a new paragraph just before the final ("DALLE, please create a
#    feature image ...") paragraph, containing a leading empty run followed
#    by a bold run with the page title text. We insert via InsertXML at a
#    collapsed range so the OOXML shape of the new paragraph matches the
#    rest of the document (an empty <w:r/> followed by the real run).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
            "<w:r/>" +
            "<w:r><w:rPr><w:b/></w:rPr><w:t>Play Apollo Rising Slot for Free - Review</w:t></w:r>" +
            "</w:p>" +
            "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'></w:p>"
$insertPoint.InsertXML($titleXml)

# InsertXML splits the content into two paragraphs and leaves a spare empty
# paragraph behind it; remove that spare paragraph so only the new title
# paragraph remains ahead of the (still unmodified) final paragraph.
$d.Paragraphs($lastIndex + 1).Range.Delete()

# 3. Replace the old DALLE image-prompt text in the final paragraph with the
#    meta-description copy, keeping the run's existing (italic) formatting.
$oldText = 'DALLE, please create a feature image for the game "Apollo Rising". The image should be in cartoon style and feature a happy Maya warrior with glasses. This should fit the space theme of the game, with the background including rich graphics of neon blue shades, stars, and rockets to create an atmosphere of a space mission. The image should capture the exciting and innovative gameplay of the game while also incorporating the Maya warrior with glasses to add a unique touch. Please make it eye-catching and attention-grabbing to draw in potential players. Thank you!'
$newText = 'Read our review of Apollo Rising slot game by IGT. Play Apollo Rising for free with 100 paylines, rocket-shaped wild symbol, and free spins bonus.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
